$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 changes
$ws.Range("G4").Value = 3.9
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 4.75
$ws.Range("Q4").Value = 2.4
$ws.Range("R4").Value = 1.53
$ws.Range("U4").Value = 2.2
$ws.Range("V4").Value = 1.62
$ws.Range("X4").Value = 19
$ws.Range("AB4").Value = 51
$ws.Range("AD4").Value = 6.5
$ws.Range("AH4").Value = 5.5
$ws.Range("AK4").Value = 17
$ws.Range("AW4").Value = 3.75

# Row 5 changes
$ws.Range("G5").Value = 1.5
$ws.Range("H5").Value = 3.9
$ws.Range("S5").Value = 1.53
$ws.Range("T5").Value = 2.38
$ws.Range("U5").Value = 2.75
$ws.Range("V5").Value = 1.4
$ws.Range("Y5").Value = 10
$ws.Range("Z5").Value = 9.5
$ws.Range("AC5").Value = 7
$ws.Range("AD5").Value = 8.5
$ws.Range("AH5").Value = 11
$ws.Range("AQ5").Value = 26
$ws.Range("AT5").Value = 2.38
$ws.Range("AW5").Value = 8.5
$ws.Range("AZ5").Value = 251

# Row 6 changes
$ws.Range("N6").Value = 8
$ws.Range("W6").Value = 11
$ws.Range("AA6").Value = 41
$ws.Range("AD6").Value = 6.5
$ws.Range("AK6").Value = 13
$ws.Range("AN6").Value = 6.5
$ws.Range("AO6").Value = 29

# Row 8 changes
$ws.Range("G8").Value = 2.55
$ws.Range("H8").Value = 3.55
$ws.Range("I8").Value = 2.4
$ws.Range("J8").Value = 3.05
$ws.Range("K8").Value = 2.25
$ws.Range("L8").Value = 2.87
$ws.Range("Q8").Value = 1.55
$ws.Range("R8").Value = 2.15
$ws.Range("U8").Value = 1.47
$ws.Range("V8").Value = 2.32
$ws.Range("X8").Value = 15
$ws.Range("Z8").Value = 29
$ws.Range("AA8").Value = 19
$ws.Range("AD8").Value = 7.3
$ws.Range("AE8").Value = 11.75
$ws.Range("AH8").Value = 11.5
$ws.Range("AL8").Value = 17
$ws.Range("AM8").Value = 21
$ws.Range("AN8").Value = 4.75
$ws.Range("AO8").Value = 13
$ws.Range("AU8").Value = 6.4
$ws.Range("BA8").Value = 65

# Row 9 changes
$ws.Range("G9").Value = 2.15
$ws.Range("I9").Value = 3.75
$ws.Range("J9").Value = 2.88
$ws.Range("L9").Value = 4.5
$ws.Range("M9").Value = 1.1
$ws.Range("N9").Value = 7
$ws.Range("Q9").Value = 2.5
$ws.Range("R9").Value = 1.5
$ws.Range("X9").Value = 9
$ws.Range("Z9").Value = 19
$ws.Range("AH9").Value = 8.5
$ws.Range("AI9").Value = 17
$ws.Range("AJ9").Value = 15
$ws.Range("AL9").Value = 41
$ws.Range("AM9").Value = 51
$ws.Range("AW9").Value = 5.5
$ws.Range("AX9").Value = 23
$ws.Range("AZ9").Value = 81
$ws.Range("BA9").Value = 126

# Row 13 changes
$ws.Range("M13").Value = 1.05
$ws.Range("P13").Value = 3.75
$ws.Range("Q13").Value = 1.93
$ws.Range("R13").Value = 1.93
